$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.369.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.20%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9973"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6272"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07491"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.80%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2904"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07738"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.846.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.000"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.69%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6808"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.16%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001058"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.10%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.25%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.106.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.67%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.188"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.424.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9989"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.487"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1375"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.422"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.99%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06511"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.416"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.478"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.096"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.093"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.833"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.17%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.142"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.96%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6960"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.70%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.581"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.265.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.837"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01836"
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "  +5.97%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9093"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.85%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9989"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.15%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.010.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -18.30%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.37"
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.740"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.59%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.079"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.98%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000117"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.29%  "

# Row 50
$ws.Range("E50").Value = "  +3.32%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.071"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
